$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.152.00'
$ws.Range('E2').Value = '  -3.57%  '
$ws.Range('D3').Value = '2.461.14'
$ws.Range('E3').Value = '  -2.76%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Value = '311.52'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('D6').Value = '93.65'
$ws.Range('E6').Value = '  -6.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.550'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.05%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '0.496'
$ws.Range('E9').Value = '  -5.11%  '
$ws.Range('D10').Value = '33.15'
$ws.Range('E10').Value = '  -7.30%  '
$ws.Range('D11').Value = '0.0776'
$ws.Range('E11').Value = '  -3.38%  '
$ws.Range('D12').Value = '0.107'
$ws.Range('E12').Value = '  -1.38%  '
$ws.Range('D13').Value = '6.94'
$ws.Range('E13').Value = '  -5.47%  '
$ws.Range('D14').Value = '2.841.10'
$ws.Range('E14').Value = '  -2.75%  '
$ws.Range('D15').Value = '2.481.17'
$ws.Range('E15').Value = '  -2.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.80'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.18%  '
$ws.Range('D17').Value = '0.781'
$ws.Range('E17').Value = '  -3.87%  '
$ws.Range('D18').Value = '41.120.94'
$ws.Range('E18').Value = '  -3.63%  '
$ws.Range('D19').Value = '6.26'
$ws.Range('E19').Value = '  -6.74%  '
$ws.Range('D20').Value = '0.0₃0918'
$ws.Range('E20').Value = '  -3.25%  '
$ws.Range('D21').Value = '11.16'
$ws.Range('E21').Value = '  -8.80%  '
$ws.Range('D22').Value = '68.27'
$ws.Range('E22').Value = '  -1.48%  '
$ws.Range('D23').Value = '235.08'
$ws.Range('E23').Value = '  -3.26%  '
$ws.Range('D24').Value = '2.74'
$ws.Range('E24').Value = '  -4.27%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('E26').Value = '  -6.27%  '
$ws.Range('D27').Value = '23.94'
$ws.Range('E27').Value = '  -5.97%  '
$ws.Range('D28').Value = '2.19'
$ws.Range('E28').Value = '  -6.08%  '
$ws.Range('D29').Value = '9.57'
$ws.Range('E29').Value = '  -5.84%  '
$ws.Range('D30').Value = '36.14'
$ws.Range('E30').Value = '  -6.24%  '
$ws.Range('D31').Value = '152.56'
$ws.Range('E31').Value = '  -4.10%  '
$ws.Range('D32').Value = '5.46'
$ws.Range('E32').Value = '  -5.26%  '
$ws.Range('E33').Value = '  -5.63%  '
$ws.Range('D34').Value = '2.56'
$ws.Range('E34').Value = '  -3.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0740'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.44%  '
$ws.Range('E36').Value = '  -4.23%  '
$ws.Range('D37').Value = '1.87'
$ws.Range('E37').Value = '  -4.77%  '
$ws.Range('D38').Value = '16.85'
$ws.Range('E38').Value = '  -8.12%  '
$ws.Range('E39').Value = '  -3.01%  '
$ws.Range('E40').Value = '  -8.46%  '
$ws.Range('D41').Value = '4.19'
$ws.Range('E41').Value = '  -1.43%  '
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('D43').Value = '19.93'
$ws.Range('E43').Value = '  -11.51%  '
$ws.Range('D44').Value = '1.975.49'
$ws.Range('E44').Value = '  -0.96%  '
$ws.Range('E45').Value = '  -5.39%  '
$ws.Range('D46').Value = '3.01'
$ws.Range('E46').Value = '  -8.41%  '
$ws.Range('E47').Value = '  -2.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '68.90'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.51%  '
$ws.Range('D49').Value = '96.47'
$ws.Range('E49').Value = '  -4.42%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.176'
$ws.Range('E50').Value = '  -6.93%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '73.71'
$ws.Range('E51').Value = '  -7.17%  '
